$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Edit 1: drop the trailing period on the "block phone numbers" bullet
#   "...provided by the user."  ->  "...provided by the user"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "The system should be able to block phone numbers provided by the user.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The system should be able to block phone numbers provided by the user",
    2
) | Out-Null

# --------------------------------------------------------------------
# Edit 2: add three new sub-bullets (3.2.1 - 3.2.2 scope items) right
# after the "Contacts Organization" bullet, matching the indentation /
# numbering (ilvl 1, numId 1) used by its sibling bullets.
# --------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$anchor.Find.Execute(
    "Contacts Organization", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null

$insertAt = $d.Range($anchor.End, $anchor.End)
$insertAt.InsertBefore(
    "`rThe system should be able to create contact groups" +
    "`rThe system should be able to add contacts to contact groups" +
    "`rThe system should be able to remove contacts from contact groups"
)

# The three freshly split paragraphs are the last three in the document;
# give them the same list level / indentation as the other sub-bullets
# (e.g. "The system should be able to add new contacts").
$total = $d.Paragraphs.Count
for ($i = $total - 2; $i -le $total; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.ListFormat.ListLevelNumber = 2
    $p.Range.ParagraphFormat.LeftIndent = 72
    $p.Range.ParagraphFormat.FirstLineIndent = -18
}
